# Apply updated TPM-derived NATMI statistics to Apoe-Vldlr LR-pair sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"27.67634766666667"
$ws.Range("H2").Value = [double]"83.029043"
$ws.Range("I2").Value = [double]"0.005965811625935536"
$ws.Range("J2").Value = [double]"0.005965811625935536"
$ws.Range("K2").Value = [double]"2"
$ws.Range("L2").Value = [double]"0.6666666666666666"
$ws.Range("M2").Value = [double]"0.1528053333333333"
$ws.Range("N2").Value = [double]"0.458416"
$ws.Range("O2").Value = [double]"0.01103433215988526"
$ws.Range("P2").Value = [double]"0.01103433215988526"
$ws.Range("Q2").Value = [double]"4.229093530654223"
$ws.Range("R2").Value = [double]"38.061841775888"
$ws.Range("S2").Value = [double]"6.582874708387786E-05"
$ws.Range("T2").Value = [double]"6.582874708387786E-05"
$ws.Range("G3").Value = [double]"27.67634766666667"
$ws.Range("H3").Value = [double]"83.029043"
$ws.Range("I3").Value = [double]"0.005965811625935536"
$ws.Range("J3").Value = [double]"0.005965811625935536"
$ws.Range("O3").Value = [double]"0.8539197603380489"
$ws.Range("P3").Value = [double]"0.8539197603380488"
$ws.Range("Q3").Value = [double]"327.2791213656014"
$ws.Range("R3").Value = [double]"2945.512092290412"
$ws.Range("S3").Value = [double]"0.005094324433840819"
$ws.Range("T3").Value = [double]"0.005094324433840818"
$ws.Range("G4").Value = [double]"27.67634766666667"
$ws.Range("H4").Value = [double]"83.029043"
$ws.Range("I4").Value = [double]"0.005965811625935536"
$ws.Range("J4").Value = [double]"0.005965811625935536"
$ws.Range("M4").Value = [double]"1.712948333333333"
$ws.Range("N4").Value = [double]"5.138845"
$ws.Range("O4").Value = [double]"0.1236949029880405"
$ws.Range("P4").Value = [double]"0.1236949029880405"
$ws.Range("Q4").Value = [double]"47.40815360837055"
$ws.Range("R4").Value = [double]"426.673382475335"
$ws.Range("S4").Value = [double]"0.0007379404903150203"
$ws.Range("T4").Value = [double]"0.0007379404903150202"
$ws.Range("G5").Value = [double]"27.67634766666667"
$ws.Range("H5").Value = [double]"83.029043"
$ws.Range("I5").Value = [double]"0.005965811625935536"
$ws.Range("J5").Value = [double]"0.005965811625935536"
$ws.Range("K5").Value = [double]"2"
$ws.Range("L5").Value = [double]"0.6666666666666666"
$ws.Range("M5").Value = [double]"0.1571906666666667"
$ws.Range("N5").Value = [double]"0.471572"
$ws.Range("O5").Value = [double]"0.01135100451402528"
$ws.Range("P5").Value = [double]"0.01135100451402528"
$ws.Range("Q5").Value = [double]"4.350463540621778"
$ws.Range("R5").Value = [double]"39.154171865596"
$ws.Range("S5").Value = [double]"6.771795469581877E-05"
$ws.Range("T5").Value = [double]"6.771795469581875E-05"
$ws.Range("I6").Value = [double]"0.009118181457976757"
$ws.Range("J6").Value = [double]"0.009118181457976757"
$ws.Range("K6").Value = [double]"2"
$ws.Range("L6").Value = [double]"0.6666666666666666"
$ws.Range("M6").Value = [double]"0.1528053333333333"
$ws.Range("N6").Value = [double]"0.458416"
$ws.Range("O6").Value = [double]"0.01103433215988526"
$ws.Range("P6").Value = [double]"0.01103433215988526"
$ws.Range("Q6").Value = [double]"6.463771341290665"
$ws.Range("R6").Value = [double]"58.173942071616"
$ws.Range("S6").Value = [double]"0.0001006130429014224"
$ws.Range("T6").Value = [double]"0.0001006130429014224"
$ws.Range("I7").Value = [double]"0.009118181457976757"
$ws.Range("J7").Value = [double]"0.009118181457976757"
$ws.Range("O7").Value = [double]"0.8539197603380489"
$ws.Range("P7").Value = [double]"0.8539197603380488"
$ws.Range("S7").Value = [double]"0.007786195325314354"
$ws.Range("T7").Value = [double]"0.007786195325314353"
$ws.Range("I8").Value = [double]"0.009118181457976757"
$ws.Range("J8").Value = [double]"0.009118181457976757"
$ws.Range("M8").Value = [double]"1.712948333333333"
$ws.Range("N8").Value = [double]"5.138845"
$ws.Range("O8").Value = [double]"0.1236949029880405"
$ws.Range("P8").Value = [double]"0.1236949029880405"
$ws.Range("Q8").Value = [double]"72.45889986024666"
$ws.Range("R8").Value = [double]"652.13009874222"
$ws.Range("S8").Value = [double]"0.001127872570871785"
$ws.Range("T8").Value = [double]"0.001127872570871785"
$ws.Range("I9").Value = [double]"0.009118181457976757"
$ws.Range("J9").Value = [double]"0.009118181457976757"
$ws.Range("K9").Value = [double]"2"
$ws.Range("L9").Value = [double]"0.6666666666666666"
$ws.Range("M9").Value = [double]"0.1571906666666667"
$ws.Range("N9").Value = [double]"0.471572"
$ws.Range("O9").Value = [double]"0.01135100451402528"
$ws.Range("P9").Value = [double]"0.01135100451402528"
$ws.Range("Q9").Value = [double]"6.649273975941333"
$ws.Range("R9").Value = [double]"59.843465783472"
$ws.Range("S9").Value = [double]"0.0001035005188891958"
$ws.Range("T9").Value = [double]"0.0001035005188891958"
$ws.Range("G10").Value = [double]"29.593002"
$ws.Range("H10").Value = [double]"88.779006"
$ws.Range("I10").Value = [double]"0.006378958578792732"
$ws.Range("J10").Value = [double]"0.006378958578792732"
$ws.Range("K10").Value = [double]"2"
$ws.Range("L10").Value = [double]"0.6666666666666666"
$ws.Range("M10").Value = [double]"0.1528053333333333"
$ws.Range("N10").Value = [double]"0.458416"
$ws.Range("O10").Value = [double]"0.01103433215988526"
$ws.Range("P10").Value = [double]"0.01103433215988526"
$ws.Range("Q10").Value = [double]"4.521968534944"
$ws.Range("R10").Value = [double]"40.697716814496"
$ws.Range("S10").Value = [double]"7.038754779254862E-05"
$ws.Range("T10").Value = [double]"7.038754779254862E-05"
$ws.Range("G11").Value = [double]"29.593002"
$ws.Range("H11").Value = [double]"88.779006"
$ws.Range("I11").Value = [double]"0.006378958578792732"
$ws.Range("J11").Value = [double]"0.006378958578792732"
$ws.Range("O11").Value = [double]"0.8539197603380489"
$ws.Range("P11").Value = [double]"0.8539197603380488"
$ws.Range("Q11").Value = [double]"349.943995854456"
$ws.Range("R11").Value = [double]"3149.495962690104"
$ws.Range("S11").Value = [double]"0.005447118780809031"
$ws.Range("T11").Value = [double]"0.00544711878080903"
$ws.Range("G12").Value = [double]"29.593002"
$ws.Range("H12").Value = [double]"88.779006"
$ws.Range("I12").Value = [double]"0.006378958578792732"
$ws.Range("J12").Value = [double]"0.006378958578792732"
$ws.Range("M12").Value = [double]"1.712948333333333"
$ws.Range("N12").Value = [double]"5.138845"
$ws.Range("O12").Value = [double]"0.1236949029880405"
$ws.Range("P12").Value = [double]"0.1236949029880405"
$ws.Range("Q12").Value = [double]"50.69128345423"
$ws.Range("R12").Value = [double]"456.22155108807"
$ws.Range("S12").Value = [double]"0.0007890446625684958"
$ws.Range("T12").Value = [double]"0.0007890446625684956"
$ws.Range("G13").Value = [double]"29.593002"
$ws.Range("H13").Value = [double]"88.779006"
$ws.Range("I13").Value = [double]"0.006378958578792732"
$ws.Range("J13").Value = [double]"0.006378958578792732"
$ws.Range("K13").Value = [double]"2"
$ws.Range("L13").Value = [double]"0.6666666666666666"
$ws.Range("M13").Value = [double]"0.1571906666666667"
$ws.Range("N13").Value = [double]"0.471572"
$ws.Range("O13").Value = [double]"0.01135100451402528"
$ws.Range("P13").Value = [double]"0.01135100451402528"
$ws.Range("Q13").Value = [double]"4.651743713048"
$ws.Range("R13").Value = [double]"41.865693417432"
$ws.Range("S13").Value = [double]"7.240758762265659E-05"
$ws.Range("T13").Value = [double]"7.240758762265657E-05"
$ws.Range("G14").Value = [double]"4539.588785666667"
$ws.Range("H14").Value = [double]"13618.766357"
$ws.Range("I14").Value = [double]"0.9785370483372949"
$ws.Range("J14").Value = [double]"0.978537048337295"
$ws.Range("K14").Value = [double]"2"
$ws.Range("L14").Value = [double]"0.6666666666666666"
$ws.Range("M14").Value = [double]"0.1528053333333333"
$ws.Range("N14").Value = [double]"0.458416"
$ws.Range("O14").Value = [double]"0.01103433215988526"
$ws.Range("P14").Value = [double]"0.01103433215988526"
$ws.Range("Q14").Value = [double]"693.6733775900568"
$ws.Range("R14").Value = [double]"6243.060398310512"
$ws.Range("S14").Value = [double]"0.01079750282210741"
$ws.Range("T14").Value = [double]"0.01079750282210741"
$ws.Range("G15").Value = [double]"4539.588785666667"
$ws.Range("H15").Value = [double]"13618.766357"
$ws.Range("I15").Value = [double]"0.9785370483372949"
$ws.Range("J15").Value = [double]"0.978537048337295"
$ws.Range("O15").Value = [double]"0.8539197603380489"
$ws.Range("P15").Value = [double]"0.8539197603380488"
$ws.Range("Q15").Value = [double]"53681.67241675147"
$ws.Range("R15").Value = [double]"483135.0517507632"
$ws.Range("S15").Value = [double]"0.8355921217980847"
$ws.Range("T15").Value = [double]"0.8355921217980847"
$ws.Range("G16").Value = [double]"4539.588785666667"
$ws.Range("H16").Value = [double]"13618.766357"
$ws.Range("I16").Value = [double]"0.9785370483372949"
$ws.Range("J16").Value = [double]"0.978537048337295"
$ws.Range("M16").Value = [double]"1.712948333333333"
$ws.Range("N16").Value = [double]"5.138845"
$ws.Range("O16").Value = [double]"0.1236949029880405"
$ws.Range("P16").Value = [double]"0.1236949029880405"
$ws.Range("Q16").Value = [double]"7776.081044426407"
$ws.Range("R16").Value = [double]"69984.72939983767"
$ws.Range("S16").Value = [double]"0.1210400452642852"
$ws.Range("T16").Value = [double]"0.1210400452642852"
$ws.Range("G17").Value = [double]"4539.588785666667"
$ws.Range("H17").Value = [double]"13618.766357"
$ws.Range("I17").Value = [double]"0.9785370483372949"
$ws.Range("J17").Value = [double]"0.978537048337295"
$ws.Range("K17").Value = [double]"2"
$ws.Range("L17").Value = [double]"0.6666666666666666"
$ws.Range("M17").Value = [double]"0.1571906666666667"
$ws.Range("N17").Value = [double]"0.471572"
$ws.Range("O17").Value = [double]"0.01135100451402528"
$ws.Range("P17").Value = [double]"0.01135100451402528"
$ws.Range("Q17").Value = [double]"713.5809876114672"
$ws.Range("R17").Value = [double]"6422.228888503204"
$ws.Range("S17").Value = [double]"0.01110737845281761"
$ws.Range("T17").Value = [double]"0.01110737845281761"
